$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GNRC")

# Row 16: Gross Margin
$ws.Range("D16").Value = 0.3727
$ws.Range("E16").Value = 0.3666
$ws.Range("F16").Value = 0.3621
$ws.Range("G16").Value = 0.3619

# Row 20: Free Cash Flow Margin
$ws.Range("D20").Value = 0.0625
$ws.Range("E20").Value = 0.0439
$ws.Range("F20").Value = -0.0007
$ws.Range("G20").Value = -0.0006

# Row 28: EBITDA Margin
$ws.Range("D28").Value = 0.1909
$ws.Range("E28").Value = 0.1765
$ws.Range("F28").Value = 0.1719
$ws.Range("G28").Value = 0.1751

# Row 29: Operating Cash Flow Margin
$ws.Range("D29").Value = 0.0644
$ws.Range("E29").Value = 0.0474
$ws.Range("F29").Value = 0.0051
$ws.Range("G29").Value = 0.0066
